$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot data (Coin, Link, Price, Volume(1h))
# Rows affected: prices/percentages updated; rows 24/25 and 49/50 swapped order.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# keeps storing them as plain text (matching the original inline-string data)
# instead of auto-converting them to numbers.

$ws.Range("D2").Value = '62.676.29'
$ws.Range("E2").Value = '  -2.35%  '
$ws.Range("D3").Value = '3.203.11'
$ws.Range("E3").Value = '  -3.49%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''592.76'
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").Value = '''136.35'
$ws.Range("E6").Value = '  -5.42%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.202.46'
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").Value = '''0.507'
$ws.Range("E9").Value = '  -3.16%  '
$ws.Range("D10").Value = '''0.145'
$ws.Range("E10").Value = '  -3.62%  '
$ws.Range("D11").Value = '''5.35'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").Value = '''0.456'
$ws.Range("E12").Value = '  -4.14%  '
$ws.Range("D13").Value = '''0.0000238'
$ws.Range("E13").Value = '  -4.66%  '
$ws.Range("D14").Value = '''33.54'
$ws.Range("E14").Value = '  -4.36%  '
$ws.Range("D15").Value = '3.730.61'
$ws.Range("E15").Value = '  -3.60%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '3.197.73'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("D18").Value = '62.745.14'
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").Value = '''6.71'
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").Value = '''463.83'
$ws.Range("E20").Value = '  -4.57%  '
$ws.Range("D21").Value = '''13.90'
$ws.Range("E21").Value = '  -3.47%  '
$ws.Range("D22").Value = '''0.715'
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").Value = '''7.68'
$ws.Range("E23").Value = '  -5.26%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''13.39'
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''84.03'
$ws.Range("E25").Value = '  -1.17%  '
$ws.Range("E27").Value = '  -3.53%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '''6.92'
$ws.Range("E29").Value = '  -4.76%  '
$ws.Range("D30").Value = '''7.86'
$ws.Range("E30").Value = '  -5.79%  '
$ws.Range("D31").Value = '''2.08'
$ws.Range("E31").Value = '  -4.51%  '
$ws.Range("D32").Value = '''27.62'
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D34").Value = '''2.43'
$ws.Range("E34").Value = '  -6.33%  '
$ws.Range("D35").Value = '''1.04'
$ws.Range("E35").Value = '  -5.48%  '
$ws.Range("D36").Value = '''5.86'
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("D37").Value = '''51.69'
$ws.Range("E37").Value = '  -3.21%  '
$ws.Range("D38").Value = '0.0₃0699'
$ws.Range("E38").Value = '  -5.59%  '
$ws.Range("E39").Value = '  -3.07%  '
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("D41").Value = '3.004.90'
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("E42").Value = '  +3.62%  '
$ws.Range("D43").Value = '''8.11'
$ws.Range("E43").Value = '  -4.63%  '
$ws.Range("D44").Value = '''2.62'
$ws.Range("E44").Value = '  -6.86%  '
$ws.Range("E45").Value = '  -6.59%  '
$ws.Range("E46").Value = '  -5.08%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").Value = '''35.55'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''25.65'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '''125.04'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("E51").Value = '  -3.04%  '
